# WS_holdings.xlsx update: refresh the "as of" disclosure date in the
# confidential banner text, and update the Weight / Percent Change figures
# for each sector row with the new model-holdings snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet ships protected; unprotect so the data cells can be written,
# then re-apply protection afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure banner (cell A16),
# in place, without disturbing the rest of the sentence.
$null = $ws.Range("A16").Replace("2021-05-28", "2021-06-09")

# Weight (column D) and Percent Change (column E) refresh for rows 2-13.
$ws.Range("D2").Value = 0.02803000413184102
$ws.Range("E2").Value = 0.002559999999999896

$ws.Range("D3").Value = 0.02206871390245104
$ws.Range("E3").Value = 0.008551068883610569

$ws.Range("D4").Value = 0.05898689401460834
$ws.Range("E4").Value = -0.0064695009242145

$ws.Range("D5").Value = 0.1363379400972748
$ws.Range("E5").Value = 0.009090909090909038

$ws.Range("D6").Value = 0.02205683209290736
$ws.Range("E6").Value = -0.005239030779305875

$ws.Range("D7").Value = 0.1268830483970565
$ws.Range("E7").Value = -0.006081337894336736

$ws.Range("D8").Value = 0.09296327786974533
$ws.Range("E8").Value = -0.01082641645615312

$ws.Range("D9").Value = 0.03186654724283192
$ws.Range("E9").Value = -0.006345957011258951

$ws.Range("D10").Value = 0.1103239994730068
$ws.Range("E10").Value = -0.01047216608487966

$ws.Range("D11").Value = 0.2813943326797173
$ws.Range("E11").Value = -0.000628817822493688

$ws.Range("D12").Value = 0.08908841009855961
$ws.Range("E12").Value = -0.001318019205422738

$ws.Range("E13").Value = -0.002427267308447445

# Restore sheet protection (original password is unknown/hashed so the
# sheet is re-protected without one).
$ws.Protect()

$wb.Save()
